# Apply QA workbook update for FY22 Medicaid load.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update run date
$ws.Range("A2").Value = "Run date: 2023-09-11"

# Widen column A to fit the new, longer row label (closest attainable
# rendering of the authored 45.71 width through the Excel width model).
$ws.Columns.Item(1).ColumnWidth = 44.8

# Insert a new row for the Fiscal Yearly Enrollment Table between the
# Yearly Enrollment Table row and the TACC server row.
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Data_warehouse Fiscal Yearly Enrollment Table"

# The distinct_members / distinct_combos figures are stored as text in
# the workbook (not numbers), so format the cells as text before typing
# the refreshed values in, across every data row (Monthly, Yearly,
# Fiscal Yearly, TACC, SPC).
$numbers = $ws.Range("B5:C9")
$numbers.NumberFormat = "@"

$ws.Range("B5").Value = "10126890"
$ws.Range("C5").Value = "436459081"
$ws.Range("B6").Value = "10126890"
$ws.Range("C6").Value = "436459081"
$ws.Range("B7").Value = "10126890"
$ws.Range("C7").Value = "436459081"
$ws.Range("B8").Value = "10126890"
$ws.Range("C8").Value = "436459081"
$ws.Range("B9").Value = "10126890"
$ws.Range("C9").Value = "436459081"

# Restore the default (non text-formatted) style footprint now that the
# values have been entered as text, keeping the sheet's appearance
# consistent with the rest of the workbook.
$numbers.ClearFormats()
